$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date column (C2:C9) from 45207 to 45208 (one day later)
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 3).Value = 45208
}

# Update the hyperlink formulas in row 2 to point to Logging_2305 instead of Logging_BRACKE
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/artfynd/A 30834-2023.xlsx", "A 30834-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/kartor/A 30834-2023.png", "A 30834-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomål/A 30834-2023.docx", "A 30834-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/klagomålsmail/A 30834-2023.docx", "A 30834-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/tillsyn/A 30834-2023.docx", "A 30834-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_2305/tillsynsmail/A 30834-2023.docx", "A 30834-2023")'
